# Mentionat explicit ca PM si GM cerute sunt valorile minime.
#
# Updates four cells on the "Schedule" sheet so the parameter descriptions
# make explicit that the phase margin / gain margin / offset voltage /
# amplification figures are the required minimum (resp. maximum) bounds,
# not just "nominal case" / "Monte Carlo analysis" values.
#
# Each cell holds rich text: a plain-formatted label run followed by a
# bold run with the actual spec value. We rewrite the full string and then
# re-apply bold to the trailing (value) portion via Characters(...).Font.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# --- A26: Margine de faza ---
$label = "Margine de faza: "
$value = "60 deg (minim, in cazul nominal)"
$ws.Range("A26").Value = $label + $value
$chars = $ws.Range("A26").Characters($label.Length + 1, $value.Length)
$chars.Font.Bold = $true
$chars.Font.Name = "Calibri"
$chars.Font.Size = 11
$chars.Font.ColorIndex = -4105

# --- A27: Margine de castig ---
$label = "Margine de castig: "
$value = "10 dB (minim, in cazul nominal)"
$ws.Range("A27").Value = $label + $value
$chars = $ws.Range("A27").Characters($label.Length + 1, $value.Length)
$chars.Font.Bold = $true
$chars.Font.Name = "Calibri"
$chars.Font.Size = 11
$chars.Font.ColorIndex = -4105

# --- A28: Amplificare ---
$label = "Amplificare: "
$value = "40 dB (minim)"
$ws.Range("A28").Value = $label + $value
$chars = $ws.Range("A28").Characters($label.Length + 1, $value.Length)
$chars.Font.Bold = $true
$chars.Font.Name = "Calibri"
$chars.Font.Size = 11
$chars.Font.ColorIndex = -4105

# --- A29: Tensiune de offset ---
$label = "Tensiune de offset: "
$value = "+/- 20 mV (maxim, in analiza Monte Carlo)"
$ws.Range("A29").Value = $label + $value
$chars = $ws.Range("A29").Characters($label.Length + 1, $value.Length)
$chars.Font.Bold = $true
$chars.Font.Name = "Calibri"
$chars.Font.Size = 11
$chars.Font.ColorIndex = -4105

# The edit session ends with the cursor on the Schedule sheet at A29 (the
# last-touched cell), matching the saved view state after this change.
$ws.Activate()
$ws.Range("A29").Select()
